$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster name labels used in columns A (Sending cluster) and D (Target cluster)
$FAPS = "FAPs"
$SCS  = "sCs"

# Ligand / Receptor symbols stay constant across all data rows
$LIGAND   = "Ccl21b"
$RECEPTOR = "Ackr4"

# Table of data rows (row number -> A, D cluster labels, then E..T numeric values)
# Columns: A(Sending cluster) B(Ligand) C(Receptor) D(Target cluster)
#          E..T numeric metrics
$rows = @(
    @{ Row=2; A=$FAPS; D=$FAPS;
       E=3; F=1; G=0.380435; H=1.141305; I=0.7997108917301441; J=0.7997108917301442;
       K=3; L=1; M=1.890486333333333; N=5.671459; O=0.9442400689667343; P=0.9442400689667344;
       Q=0.7192071682216667; R=6.472864513995001; S=0.7551190675607199; T=0.75511906756072 },

    @{ Row=3; A=$FAPS; D=$SCS;
       E=3; F=1; G=0.380435; H=1.141305; I=0.7997108917301441; J=0.7997108917301442;
       K=2; L=0.6666666666666666; M=0.1116383333333333; N=0.334915; O=0.05575993103326565; P=0.05575993103326566;
       Q=0.04247112934166667; R=0.382240164075; S=0.04459182416942421; T=0.04459182416942422 },

    @{ Row=4; A=$SCS; D=$FAPS;
       E=2; F=0.6666666666666666; G=0.09528066666666667; H=0.285842; I=0.2002891082698559; J=0.2002891082698559;
       K=3; L=1; M=1.890486333333333; N=5.671459; O=0.9442400689667343; P=0.9442400689667344;
       Q=0.1801267981642222; R=1.621141183478; S=0.1891210014060144; T=0.1891210014060145 },

    @{ Row=5; A=$SCS; D=$SCS;
       E=2; F=0.6666666666666666; G=0.09528066666666667; H=0.285842; I=0.2002891082698559; J=0.2002891082698559;
       K=2; L=0.6666666666666666; M=0.1116383333333333; N=0.334915; O=0.05575993103326565; P=0.05575993103326566;
       Q=0.01063697482555556; R=0.09573277343; S=0.01116810686384144; T=0.01116810686384144 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value  = $r.A
    $ws.Cells.Item($rowNum, 2).Value  = $LIGAND
    $ws.Cells.Item($rowNum, 3).Value  = $RECEPTOR
    $ws.Cells.Item($rowNum, 4).Value  = $r.D

    $ws.Cells.Item($rowNum, 5).Value  = $r.E
    $ws.Cells.Item($rowNum, 6).Value  = $r.F
    $ws.Cells.Item($rowNum, 7).Value  = $r.G
    $ws.Cells.Item($rowNum, 8).Value  = $r.H
    $ws.Cells.Item($rowNum, 9).Value  = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
}
